$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (Fri Aug 18 17:28:40 UTC 2023 GitHub Actions run).
# Rows 23/24 (Chainlink <-> BinanceUSD) and rows 29/30 (Hedera <-> Toncoin) swapped
# position in the refreshed ranking, so Coin/Link/Price/Volume are updated together.
# Price column (D) values are plain-text (not numbers) in the source data, e.g.
# "25.979.81" or "1.009" -- prefixing with an apostrophe forces Excel to keep them
# as text instead of auto-converting to a number (which would also drop trailing
# zeros like "1.010" -> 1.01). ClearFormats() afterwards drops the quote-prefix
# style so no stray number-format style is left on the cell.
$updates = @(
    @{ Cell = 'D2'; Value = '25.947.29'; ForceText = $true },
    @{ Cell = 'E2'; Value = '  -7.47%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.660.93'; ForceText = $true },
    @{ Cell = 'E3'; Value = '  -4.83%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.009'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.74%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '217.67'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -4.03%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.5012'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -13.57%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '1.009'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  +0.63%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.2621'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -3.44%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.06292'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -4.70%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '21.32'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -7.80%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.07359'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -2.05%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '1.665.41'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -4.39%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '4.528'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -4.24%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '0.5713'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -5.56%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '1.894.50'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -4.46%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '0.000008350'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -3.60%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '64.23'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -13.48%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '26.024.39'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -7.19%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '4.911'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -7.87%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '1.008'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +0.65%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '10.71'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -4.97%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '185.52'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -9.55%  '; ForceText = $false },
    @{ Cell = 'B23'; Value = 'BinanceUSD'; ForceText = $false },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false },
    @{ Cell = 'D23'; Value = '1.010'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +0.74%  '; ForceText = $false },
    @{ Cell = 'B24'; Value = 'Chainlink'; ForceText = $false },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; ForceText = $false },
    @{ Cell = 'D24'; Value = '6.143'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -7.37%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '142.57'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -4.87%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '7.609'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -5.24%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '0.1163'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -5.70%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '15.64'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -3.11%  '; ForceText = $false },
    @{ Cell = 'B29'; Value = 'Toncoin'; ForceText = $false },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false },
    @{ Cell = 'D29'; Value = '1.299'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -6.23%  '; ForceText = $false },
    @{ Cell = 'B30'; Value = 'Hedera'; ForceText = $false },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false },
    @{ Cell = 'D30'; Value = '0.05801'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -6.20%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '1.318'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -5.23%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '3.477'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -6.91%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '3.478'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -6.41%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '1.636'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -2.45%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.9979'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -3.67%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.5950'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -6.68%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '2.368'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -3.85%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '2.634'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -3.08%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.01590'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -5.08%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '1.077.00'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -4.46%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '5.936'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -4.53%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '0.8531'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -2.30%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.008'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +0.49%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '99.28'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -0.32%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.816.16'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +2.55%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '55.59'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -6.51%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '1.005'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +0.72%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '8.041'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -2.24%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.4305'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -2.61%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.05171'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -3.90%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.Value = "'" + $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
